$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts D:L data right (old D..K -> E..L)
$ws.Range("D:D").Insert()

# Copy format (number format, font, etc.) from column E into column D so the new column matches
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)  # xlPasteFormats

$v = $ws.Range("D7").NumberFormat
Write-Host "D7 NumberFormat=$v"
$v2 = $ws.Range("D8").NumberFormat
Write-Host "D8 NumberFormat=$v2"
